# Update the "想去人数" (interest count) figures that were refreshed when the
# gh-pages data was regenerated (commit: "Update gh-pages to output generated
# at 456a3b4").
#
# Sheet "展览" (Exhibitions) - column F, rows 2-11
# Sheet "演出" (Performances) - column F, row 3
# Sheet "全部类型" (All types) - column F, rows 2-13 (union of the two sheets above)

$wb = $excel.ActiveWorkbook

$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 7665
$wsExpo.Range("F3").Value = 295
$wsExpo.Range("F4").Value = 31
$wsExpo.Range("F5").Value = 468
$wsExpo.Range("F6").Value = 4344
$wsExpo.Range("F8").Value = 610
$wsExpo.Range("F9").Value = 280
$wsExpo.Range("F10").Value = 687
$wsExpo.Range("F11").Value = 165

$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F3").Value = 13

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 7665
$wsAll.Range("F4").Value = 295
$wsAll.Range("F5").Value = 31
$wsAll.Range("F6").Value = 468
$wsAll.Range("F7").Value = 4344
$wsAll.Range("F9").Value = 610
$wsAll.Range("F10").Value = 280
$wsAll.Range("F11").Value = 687
$wsAll.Range("F12").Value = 13
$wsAll.Range("F13").Value = 165
